# Update notebook for 10-30-2020
# Append the new daily data row (row 96) to Sheet1:
#   A96 = 2020-10-29 (serial 44133), B96 = 7170 (Total Cases),
#   C96 = New Cases (B96-B95), D96 = 7-Day Average (AVERAGE(C90:C96))

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New data row ---
$ws.Range("A96").Value = 44133
$ws.Range("B96").Value = 7170
$ws.Range("C96").Formula = "=B96-B95"
$ws.Range("D96").Formula = "=AVERAGE(C90:C96)"

# Match the "New Cases" column's integer number format used by the rest
# of column C (column A/D already inherit their formatting automatically).
$ws.Range("C96").NumberFormat = "0"

# --- Update the active selection to follow the newly-entered row,
#     matching Excel's normal behaviour of moving to the next empty row ---
$ws.Activate() | Out-Null
$ws.Range("D97").Select() | Out-Null
